$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before row 21; this shifts the existing Stagecoach..AEC
# Routemaster block (rows 21-26) down to rows 23-28.
$ws.Rows("21:22").Insert()

# Append a new row 29 (BMMO CM5T) at the end of the table.
# Set the Vehicle name first so the "BMMO CM5T" string becomes shared string
# index 42 (it is added to the shared string table before the ERF entries).
$ws.Range("A29").Value = "BMMO CM5T"
$ws.Range("B29").Value = 1958
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "Bus"
$ws.Range("E29").Formula = "=IF(B29 > 1900, ((B29-1900)*10)+400+C29, ((B29-1730)*2)+C29)+VLOOKUP(D29,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Range("F29").Value = 76
$ws.Range("G29").Value = 34
$ws.Range("H29").Formula = "=SQRT(F29*G29)/`$B`$1"
$ws.Range("H29").NumberFormat = "0"
$ws.Range("I29").Formula = "=H29*0.9"
$ws.Range("I29").NumberFormat = "0"
$ws.Range("J29").Value = "x"
$ws.Range("J29").NumberFormat = "0"

# Row 21: ERF C-Series (Heavy Goods)
$ws.Range("A21").Value = "ERF C-Series"
$ws.Range("B21").Value = 1982
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "Heavy Goods"
$ws.Range("E21").Formula = "=IF(B21 > 1900, ((B21-1900)*10)+400+C21, ((B21-1730)*2)+C21)+VLOOKUP(D21,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Range("F21").Value = 60
$ws.Range("G21").Value = 28
$ws.Range("H21").Formula = "=SQRT(F21*G21)/`$B`$1"
$ws.Range("I21").Formula = "=H21*0.9"
$ws.Range("J21").Clear()

# Row 22: ERF E-Series (Heavy Goods)
$ws.Range("A22").Value = "ERF E-Series"
$ws.Range("B22").Value = 1986
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "Heavy Goods"
$ws.Range("E22").Formula = "=IF(B22 > 1900, ((B22-1900)*10)+400+C22, ((B22-1730)*2)+C22)+VLOOKUP(D22,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Range("F22").Value = 62
$ws.Range("G22").Value = 30
$ws.Range("H22").Formula = "=SQRT(F22*G22)/`$B`$1"
$ws.Range("I22").Formula = "=H22*0.9"
$ws.Range("J22").Clear()

# Update the view: scroll down a bit and select G21, matching the saved file.
[void]$ws.Range("G21").Select()

Write-Host ("UsedRange: " + $ws.UsedRange.Address())
